$d = $word.ActiveDocument

# 1. Remove the stale _GoBack bookmark from its old location
#    (between "... on " and "the project name ...")
try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
} catch {
    # bookmark not present; nothing to remove
}

# 2. Insert a new list item right after the "Multi-threaded Debug DLL"
#    paragraph and before the "Now it should work" paragraph, and move
#    the _GoBack bookmark to the end of this new paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Multi-threaded Debug DLL*") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Run </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>using</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> x86 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Debug</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>depends</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>your</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Visual Studio most </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>likely</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$newPara.Range.InsertXML($newParaXml) | Out-Null
